# patterns/transaction-complex-read.pptx
# commit: "fix part of #53 with some TBDs"
#
# 1) Notes Master date placeholder: cached datetimeFigureOut text
#    2023/1/12 -> 2023/2/23
# 2) Slide 1, shape "矩形 8" (RESULT / SQL annotation callout box):
#    remove the whole "COUNT(DISTINCT medium)," paragraph (2nd of 3),
#    leaving "RESULT" and "medium.type, medium.id" paragraphs intact.

$p = $ppt.ActivePresentation

# --- 1) Notes Master date field -------------------------------------------
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $nmShape = $nm.Shapes.Item($i)
    if ($nmShape.HasTextFrame -and $nmShape.TextFrame.HasText) {
        if ($nmShape.TextFrame.TextRange.Text.TrimEnd() -eq "2023/1/12") {
            $nmShape.TextFrame.TextRange.Text = "2023/2/23"
        }
    }
}

# --- 2) Slide 1 "矩形 8" shape: drop the COUNT(DISTINCT medium), line -----
$slide = $p.Slides.Item(1)
$outerGroup = $slide.Shapes.Item(1)

$resultBox = $null
for ($i = 1; $i -le $outerGroup.GroupItems.Count; $i++) {
    $candidate = $outerGroup.GroupItems.Item($i)
    if ($candidate.Id -eq 9) {
        $resultBox = $candidate
    }
}

$tr = $resultBox.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.TrimEnd() -eq "COUNT(DISTINCT medium),") {
        $para.Delete()
    }
}
